$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPRINT BACKLOG 1")

# --- Row 11: "Nro de Horas" label + SUM formulas (creates the burndown series) ---
$ws.Range("B11").Value = "Nro de Horas"
$ws.Range("C11:I11").Formula = "=SUM(C6:C10)"

# --- Row 25: author added to the authors list ---
$ws.Range("A25").Value = "Carlos Zarate Carpio"

# --- External reference to another workbook (Hoja1) ---
# Writing then clearing a formula that references the external workbook
# registers the external link part without leaving a stray value behind.
$tmpCell = $ws.Range("Z1000")
$tmpCell.Formula = "='[Hoja1.xlsx]Hoja1'!A1"
$tmpCell.ClearContents()

# --- Sprint Burndown chart ---
$co = $ws.ChartObjects().Add(228600, 2628900, 4305300, 2628900)
$chart = $co.Chart
$chart.ChartType = 4   # xlLine

$ser = $chart.SeriesCollection().NewSeries()
# Build the series via an explicit SERIES() formula (quoted, absolute refs) so
# the sheet name with spaces is preserved correctly in the saved chart XML.
$ser.Formula = "=SERIES('SPRINT BACKLOG 1'!`$B`$11,'SPRINT BACKLOG 1'!`$C`$4:`$I`$4,'SPRINT BACKLOG 1'!`$C`$11:`$I`$11,1)"

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Sprint Burndown 1"
$chart.HasLegend = $false

$chart.Axes(1).HasTitle = $true
$chart.Axes(1).AxisTitle.Text = "Nro Dias"
$chart.Axes(2).HasTitle = $true
$chart.Axes(2).AxisTitle.Text = "Nro de Horas"

# --- Selection as left by the author ---
$ws.Range("B21").Select() | Out-Null

Write-Host "edit applied"
